$wb = $excel.ActiveWorkbook

# --- Step 1: update regression coefficient/p-value results (rerun models) ---
$ws = $wb.Worksheets.Item("summ9")
$ws.Cells.Item(2, 2).Value = 2.91977941239301
$ws.Cells.Item(2, 3).Value = 0.02020854018698176
$ws.Cells.Item(3, 2).Value = -0.8836525202627834
$ws.Cells.Item(3, 3).Value = 0.2511759660377603
$ws.Cells.Item(4, 2).Value = -1.658292534450327
$ws.Cells.Item(4, 3).Value = 0.0001950766152178553
$ws.Cells.Item(5, 2).Value = -1.269960975940354
$ws.Cells.Item(5, 3).Value = 0.005465805644036517
$ws.Cells.Item(6, 2).Value = -1.717955301327796
$ws.Cells.Item(6, 3).Value = 0.001773356918846869
$ws.Cells.Item(7, 2).Value = 0.02063679444222259
$ws.Cells.Item(7, 3).Value = 0.9485362906088662
$ws.Cells.Item(8, 2).Value = 0.0009722337068753616
$ws.Cells.Item(8, 3).Value = 0.000000007048020594797252
$ws.Cells.Item(9, 2).Value = -0.01529143487216724
$ws.Cells.Item(9, 3).Value = 0.02814576229057139
$ws.Cells.Item(10, 2).Value = 0.9374037314912899
$ws.Cells.Item(10, 3).Value = 0.0006257033175050727
$ws.Cells.Item(11, 2).Value = 1.007664095354128
$ws.Cells.Item(11, 3).Value = 0.0004973160060871666
$ws.Cells.Item(12, 2).Value = 0.3406214698816832
$ws.Cells.Item(12, 3).Value = 0.3236299926925081
$ws.Cells.Item(13, 2).Value = -0.0000426881883797294
$ws.Cells.Item(13, 3).Value = 0.6936836570081111
$ws.Cells.Item(14, 2).Value = -0.0000001507768201921867
$ws.Cells.Item(14, 3).Value = 0.0904614521892352
$ws.Cells.Item(15, 2).Value = -0.09356027058660131
$ws.Cells.Item(15, 3).Value = 0.3829752554260589
$ws.Cells.Item(16, 2).Value = 0.1377436362800698
$ws.Cells.Item(16, 3).Value = 0.02294691547305563
$ws.Cells.Item(17, 2).Value = -5.12196525336998
$ws.Cells.Item(17, 3).Value = 0.4958629284127962
$ws.Cells.Item(18, 2).Value = -0.01985878809752055
$ws.Cells.Item(18, 3).Value = 0.03378864884664836
$ws.Cells.Item(19, 2).Value = -0.007309094848754527
$ws.Cells.Item(19, 3).Value = 0.227564989141302
$ws.Cells.Item(20, 2).Value = 0.1303250816089186
$ws.Cells.Item(20, 3).Value = 0.9021139803612203
$ws.Cells.Item(21, 2).Value = 1.600557203597112
$ws.Cells.Item(21, 3).Value = 0.4437114345029356
$ws.Cells.Item(22, 2).Value = 0.0006277950442059815
$ws.Cells.Item(22, 3).Value = 0.01885601983404929

$ws = $wb.Worksheets.Item("summ7")
$ws.Cells.Item(2, 2).Value = 2.965739131530457
$ws.Cells.Item(2, 3).Value = 0.01895014120932965
$ws.Cells.Item(3, 2).Value = -0.2801659687415298
$ws.Cells.Item(3, 3).Value = 0.7367477483065054
$ws.Cells.Item(4, 2).Value = -1.938515504509885
$ws.Cells.Item(4, 3).Value = 0.00001951193049142829
$ws.Cells.Item(5, 2).Value = -1.677374849232944
$ws.Cells.Item(5, 3).Value = 0.0002883801997917914
$ws.Cells.Item(6, 2).Value = -1.711805216210918
$ws.Cells.Item(6, 3).Value = 0.001465065438236873
$ws.Cells.Item(7, 2).Value = -0.210961603071576
$ws.Cells.Item(7, 3).Value = 0.4966723216900156
$ws.Cells.Item(8, 2).Value = 0.001051952049913179
$ws.Cells.Item(8, 3).Value = 0.000000001842988830777967
$ws.Cells.Item(9, 2).Value = -0.01523942611568955
$ws.Cells.Item(9, 3).Value = 0.03251725927155604
$ws.Cells.Item(10, 2).Value = 0.8509577709608712
$ws.Cells.Item(10, 3).Value = 0.00235849425649131
$ws.Cells.Item(11, 2).Value = 1.005539855621983
$ws.Cells.Item(11, 3).Value = 0.0007781230081573885
$ws.Cells.Item(12, 2).Value = 0.3172386620128906
$ws.Cells.Item(12, 3).Value = 0.3687187912908578
$ws.Cells.Item(13, 2).Value = -0.00009798893760572084
$ws.Cells.Item(13, 3).Value = 0.3757796934307363
$ws.Cells.Item(14, 2).Value = -0.0000001180133366427938
$ws.Cells.Item(14, 3).Value = 0.1958088061438568
$ws.Cells.Item(15, 2).Value = -0.1233827631192954
$ws.Cells.Item(15, 3).Value = 0.2772227494393097
$ws.Cells.Item(16, 2).Value = 0.1747566430474503
$ws.Cells.Item(16, 3).Value = 0.00846347458408438
$ws.Cells.Item(17, 2).Value = 3.143601408456903
$ws.Cells.Item(17, 3).Value = 0.68540218608651
$ws.Cells.Item(18, 2).Value = -0.01682225949514497
$ws.Cells.Item(18, 3).Value = 0.07015596357505013
$ws.Cells.Item(19, 2).Value = -0.00589458755546486
$ws.Cells.Item(19, 3).Value = 0.3175928449331434
$ws.Cells.Item(20, 2).Value = 0.3658199029992339
$ws.Cells.Item(20, 3).Value = 0.7414220893180021
$ws.Cells.Item(21, 2).Value = 0.3270611896776602
$ws.Cells.Item(21, 3).Value = 0.8762591635188152
$ws.Cells.Item(22, 2).Value = 0.0006512567742717741
$ws.Cells.Item(22, 3).Value = 0.01494274443989387

$ws = $wb.Worksheets.Item("summ13")
$ws.Cells.Item(2, 2).Value = 1.991099410031772
$ws.Cells.Item(2, 3).Value = 0.1030671730356838
$ws.Cells.Item(3, 2).Value = -0.9415686714075823
$ws.Cells.Item(3, 3).Value = 0.1851245578519826
$ws.Cells.Item(4, 2).Value = -1.5982895977745
$ws.Cells.Item(4, 3).Value = 0.0001627454558476533
$ws.Cells.Item(5, 2).Value = -1.282115882534295
$ws.Cells.Item(5, 3).Value = 0.003343740408029857
$ws.Cells.Item(6, 2).Value = -1.609586112084
$ws.Cells.Item(6, 3).Value = 0.002233339035562454
$ws.Cells.Item(7, 2).Value = 0.02602275518398819
$ws.Cells.Item(7, 3).Value = 0.9282900051189252
$ws.Cells.Item(8, 2).Value = 0.001002170470295079
$ws.Cells.Item(8, 3).Value = 0.000000004533323404310383
$ws.Cells.Item(9, 2).Value = -0.01677484079286197
$ws.Cells.Item(9, 3).Value = 0.01734722379882828
$ws.Cells.Item(10, 2).Value = 0.9845615752870482
$ws.Cells.Item(10, 3).Value = 0.0004052682110514608
$ws.Cells.Item(11, 2).Value = 0.9988100810614908
$ws.Cells.Item(11, 3).Value = 0.0006019163985315646
$ws.Cells.Item(12, 2).Value = 0.5035359424879496
$ws.Cells.Item(12, 3).Value = 0.1597806288220962
$ws.Cells.Item(13, 2).Value = -0.0001317468710046529
$ws.Cells.Item(13, 3).Value = 0.2247259677032915
$ws.Cells.Item(14, 2).Value = -0.00000006525905997098933
$ws.Cells.Item(14, 3).Value = 0.465763567957579
$ws.Cells.Item(15, 2).Value = -0.1457790323131748
$ws.Cells.Item(15, 3).Value = 0.1685777904653972
$ws.Cells.Item(16, 2).Value = 0.1359066887561832
$ws.Cells.Item(16, 3).Value = 0.02712221451955548
$ws.Cells.Item(17, 2).Value = 1.288951086713106
$ws.Cells.Item(17, 3).Value = 0.8643746905932619
$ws.Cells.Item(18, 2).Value = -0.009929016963739122
$ws.Cells.Item(18, 3).Value = 0.281894446285663
$ws.Cells.Item(19, 2).Value = -0.0004793725200102507
$ws.Cells.Item(19, 3).Value = 0.9345642508819416
$ws.Cells.Item(20, 2).Value = 0.2226533188100894
$ws.Cells.Item(20, 3).Value = 0.8306123067552318
$ws.Cells.Item(21, 2).Value = -0.6887485282733223
$ws.Cells.Item(21, 3).Value = 0.7353977647407005
$ws.Cells.Item(22, 2).Value = 0.0005351031689666172
$ws.Cells.Item(22, 3).Value = 0.04323327003757856

$ws = $wb.Worksheets.Item("summ10")
$ws.Cells.Item(2, 2).Value = 2.144578677633742
$ws.Cells.Item(2, 3).Value = 0.08246325822970688
$ws.Cells.Item(3, 2).Value = -1.005517352620882
$ws.Cells.Item(3, 3).Value = 0.1642294486989858
$ws.Cells.Item(4, 2).Value = -1.445740172834051
$ws.Cells.Item(4, 3).Value = 0.0006287751820513697
$ws.Cells.Item(5, 2).Value = -1.144171541408449
$ws.Cells.Item(5, 3).Value = 0.008448436071828817
$ws.Cells.Item(6, 2).Value = -1.792156360184232
$ws.Cells.Item(6, 3).Value = 0.0005970813962540921
$ws.Cells.Item(7, 2).Value = 0.07008097669295234
$ws.Cells.Item(7, 3).Value = 0.8107386487175573
$ws.Cells.Item(8, 2).Value = 0.0011353654518695
$ws.Cells.Item(8, 3).Value = 0.00000000005873431967648877
$ws.Cells.Item(9, 2).Value = -0.0185597138258779
$ws.Cells.Item(9, 3).Value = 0.007446198832485831
$ws.Cells.Item(10, 2).Value = 0.7509996971399185
$ws.Cells.Item(10, 3).Value = 0.006171121368707644
$ws.Cells.Item(11, 2).Value = 0.9218135420783514
$ws.Cells.Item(11, 3).Value = 0.001658020823129727
$ws.Cells.Item(12, 2).Value = 0.4320035709343067
$ws.Cells.Item(12, 3).Value = 0.2032910423676153
$ws.Cells.Item(13, 2).Value = -0.0001097669273750651
$ws.Cells.Item(13, 3).Value = 0.3147942751206015
$ws.Cells.Item(14, 2).Value = -0.00000004981509231007353
$ws.Cells.Item(14, 3).Value = 0.5770543635068671
$ws.Cells.Item(15, 2).Value = -0.1691491386476577
$ws.Cells.Item(15, 3).Value = 0.1163613046239695
$ws.Cells.Item(16, 2).Value = 0.1305244206655957
$ws.Cells.Item(16, 3).Value = 0.03880790125441919
$ws.Cells.Item(17, 2).Value = 2.009748158815011
$ws.Cells.Item(17, 3).Value = 0.789605335621824
$ws.Cells.Item(18, 2).Value = -0.01758682119042814
$ws.Cells.Item(18, 3).Value = 0.05760535745457938
$ws.Cells.Item(19, 2).Value = 0.0004080014855090707
$ws.Cells.Item(19, 3).Value = 0.9472394252025537
$ws.Cells.Item(20, 2).Value = 0.5350054779968398
$ws.Cells.Item(20, 3).Value = 0.6164879016959272
$ws.Cells.Item(21, 2).Value = -0.4383191424580066
$ws.Cells.Item(21, 3).Value = 0.8299435935455477
$ws.Cells.Item(22, 2).Value = 0.0004819889675413649
$ws.Cells.Item(22, 3).Value = 0.06933436627592109

$ws = $wb.Worksheets.Item("summ3")
$ws.Cells.Item(2, 2).Value = 2.899596024877852
$ws.Cells.Item(2, 3).Value = 0.01908978330257155
$ws.Cells.Item(3, 2).Value = -0.9668907970099783
$ws.Cells.Item(3, 3).Value = 0.1791849963956914
$ws.Cells.Item(4, 2).Value = -1.578631227114108
$ws.Cells.Item(4, 3).Value = 0.0002585181669944186
$ws.Cells.Item(5, 2).Value = -1.24301350927264
$ws.Cells.Item(5, 3).Value = 0.005158713455261608
$ws.Cells.Item(6, 2).Value = -1.876082499478784
$ws.Cells.Item(6, 3).Value = 0.0003533873097970665
$ws.Cells.Item(7, 2).Value = -0.01977534083577642
$ws.Cells.Item(7, 3).Value = 0.94703642329268
$ws.Cells.Item(8, 2).Value = 0.0010614921530932
$ws.Cells.Item(8, 3).Value = 0.0000000004591183418294408
$ws.Cells.Item(9, 2).Value = -0.01745283197371262
$ws.Cells.Item(9, 3).Value = 0.01131410630153422
$ws.Cells.Item(10, 2).Value = 0.8318294441446272
$ws.Cells.Item(10, 3).Value = 0.002301144246287924
$ws.Cells.Item(11, 2).Value = 1.038254379077723
$ws.Cells.Item(11, 3).Value = 0.0003168352845008396
$ws.Cells.Item(12, 2).Value = 0.4092601391874487
$ws.Cells.Item(12, 3).Value = 0.2290010446527257
$ws.Cells.Item(13, 2).Value = -0.0001058171645634919
$ws.Cells.Item(13, 3).Value = 0.3394510197173222
$ws.Cells.Item(14, 2).Value = -0.00000006722298300313112
$ws.Cells.Item(14, 3).Value = 0.4610585188301545
$ws.Cells.Item(15, 2).Value = -0.1642656629235547
$ws.Cells.Item(15, 3).Value = 0.1220346966993756
$ws.Cells.Item(16, 2).Value = 0.1360111405698775
$ws.Cells.Item(16, 3).Value = 0.02676118366451358
$ws.Cells.Item(17, 2).Value = -2.404035767257663
$ws.Cells.Item(17, 3).Value = 0.7526328220330172
$ws.Cells.Item(18, 2).Value = -0.01329209193273824
$ws.Cells.Item(18, 3).Value = 0.1514433861968509
$ws.Cells.Item(19, 2).Value = -0.004288864222076437
$ws.Cells.Item(19, 3).Value = 0.4688408910202122
$ws.Cells.Item(20, 2).Value = -0.1832855092909823
$ws.Cells.Item(20, 3).Value = 0.8647906241186917
$ws.Cells.Item(21, 2).Value = -1.089918868547423
$ws.Cells.Item(21, 3).Value = 0.603959545870117
$ws.Cells.Item(22, 2).Value = 0.0005729725585306499
$ws.Cells.Item(22, 3).Value = 0.03103930110799209

$ws = $wb.Worksheets.Item("summ5")
$ws.Cells.Item(2, 2).Value = 2.178095285652875
$ws.Cells.Item(2, 3).Value = 0.07735898807490632
$ws.Cells.Item(3, 2).Value = -0.8438760267191386
$ws.Cells.Item(3, 3).Value = 0.2478992780623479
$ws.Cells.Item(4, 2).Value = -1.617432611632353
$ws.Cells.Item(4, 3).Value = 0.000142846394768782
$ws.Cells.Item(5, 2).Value = -1.320489307418281
$ws.Cells.Item(5, 3).Value = 0.00249748185259218
$ws.Cells.Item(6, 2).Value = -1.407092690603103
$ws.Cells.Item(6, 3).Value = 0.009121242931622059
$ws.Cells.Item(7, 2).Value = -0.04611546055016356
$ws.Cells.Item(7, 3).Value = 0.8755418424837209
$ws.Cells.Item(8, 2).Value = 0.001092409806010313
$ws.Cells.Item(8, 3).Value = 0.000000000280333104236044
$ws.Cells.Item(9, 2).Value = -0.01267132120471529
$ws.Cells.Item(9, 3).Value = 0.06961765179096918
$ws.Cells.Item(10, 2).Value = 0.8850529159647539
$ws.Cells.Item(10, 3).Value = 0.00123912007395106
$ws.Cells.Item(11, 2).Value = 0.9475172918541884
$ws.Cells.Item(11, 3).Value = 0.001239751820901147
$ws.Cells.Item(12, 2).Value = 0.2689522924316408
$ws.Cells.Item(12, 3).Value = 0.4409283710964089
$ws.Cells.Item(13, 2).Value = -0.00003783128238642764
$ws.Cells.Item(13, 3).Value = 0.7307929240138296
$ws.Cells.Item(14, 2).Value = -0.0000001474008447574142
$ws.Cells.Item(14, 3).Value = 0.1021178275945915
$ws.Cells.Item(15, 2).Value = -0.1674716155531779
$ws.Cells.Item(15, 3).Value = 0.1185868891938509
$ws.Cells.Item(16, 2).Value = 0.1090963763761736
$ws.Cells.Item(16, 3).Value = 0.07912441899495995
$ws.Cells.Item(17, 2).Value = 1.781542328496833
$ws.Cells.Item(17, 3).Value = 0.8197716645113906
$ws.Cells.Item(18, 2).Value = -0.01423905384419162
$ws.Cells.Item(18, 3).Value = 0.1354629069789995
$ws.Cells.Item(19, 2).Value = 0.0002325487004053381
$ws.Cells.Item(19, 3).Value = 0.9700652129736521
$ws.Cells.Item(20, 2).Value = -0.0972854882970741
$ws.Cells.Item(20, 3).Value = 0.9261070372421568
$ws.Cells.Item(21, 2).Value = 1.326490313592177
$ws.Cells.Item(21, 3).Value = 0.5329765428024249
$ws.Cells.Item(22, 2).Value = 0.0003898571398815965
$ws.Cells.Item(22, 3).Value = 0.1531947754767564

$ws = $wb.Worksheets.Item("summ14")
$ws.Cells.Item(2, 2).Value = 2.470294074890139
$ws.Cells.Item(2, 3).Value = 0.04876724632254324
$ws.Cells.Item(3, 2).Value = -1.238222341481711
$ws.Cells.Item(3, 3).Value = 0.09327610118048969
$ws.Cells.Item(4, 2).Value = -1.573613845289684
$ws.Cells.Item(4, 3).Value = 0.000439660055782577
$ws.Cells.Item(5, 2).Value = -1.330976344721692
$ws.Cells.Item(5, 3).Value = 0.003639102336446934
$ws.Cells.Item(6, 2).Value = -1.742869302602118
$ws.Cells.Item(6, 3).Value = 0.001603536804519054
$ws.Cells.Item(7, 2).Value = 0.1361153591612613
$ws.Cells.Item(7, 3).Value = 0.6604640440224586
$ws.Cells.Item(8, 2).Value = 0.001046090900104195
$ws.Cells.Item(8, 3).Value = 0.000000001446867753295286
$ws.Cells.Item(9, 2).Value = -0.01847154671797234
$ws.Cells.Item(9, 3).Value = 0.01007309098870214
$ws.Cells.Item(10, 2).Value = 0.8091585234670733
$ws.Cells.Item(10, 3).Value = 0.002817313413413157
$ws.Cells.Item(11, 2).Value = 0.9973903734068877
$ws.Cells.Item(11, 3).Value = 0.0006862841958236924
$ws.Cells.Item(12, 2).Value = 0.4366616510440364
$ws.Cells.Item(12, 3).Value = 0.2157095977001076
$ws.Cells.Item(13, 2).Value = -0.00006280644343798359
$ws.Cells.Item(13, 3).Value = 0.5576520992588134
$ws.Cells.Item(14, 2).Value = -0.0000001395438850824026
$ws.Cells.Item(14, 3).Value = 0.1141186960794569
$ws.Cells.Item(15, 2).Value = -0.1081518362070578
$ws.Cells.Item(15, 3).Value = 0.3078399125891539
$ws.Cells.Item(16, 2).Value = 0.1021030127256359
$ws.Cells.Item(16, 3).Value = 0.10093695843114
$ws.Cells.Item(17, 2).Value = -1.26921718922784
$ws.Cells.Item(17, 3).Value = 0.8645337878808645
$ws.Cells.Item(18, 2).Value = -0.01833742511010925
$ws.Cells.Item(18, 3).Value = 0.0474253518615123
$ws.Cells.Item(19, 2).Value = -0.00293492170552592
$ws.Cells.Item(19, 3).Value = 0.6172026334469802
$ws.Cells.Item(20, 2).Value = 0.306234231150161
$ws.Cells.Item(20, 3).Value = 0.7718140068483554
$ws.Cells.Item(21, 2).Value = 0.9689758200931483
$ws.Cells.Item(21, 3).Value = 0.6311630388898617
$ws.Cells.Item(22, 2).Value = 0.000617745876504495
$ws.Cells.Item(22, 3).Value = 0.01878296656451707

$ws = $wb.Worksheets.Item("summ1")
$ws.Cells.Item(2, 2).Value = 1.926320850408461
$ws.Cells.Item(2, 3).Value = 0.1124978078024224
$ws.Cells.Item(3, 2).Value = -0.8979156853737346
$ws.Cells.Item(3, 3).Value = 0.2121559583743522
$ws.Cells.Item(4, 2).Value = -1.497911516258987
$ws.Cells.Item(4, 3).Value = 0.0003629929901182239
$ws.Cells.Item(5, 2).Value = -1.063749643348088
$ws.Cells.Item(5, 3).Value = 0.01449513765571473
$ws.Cells.Item(6, 2).Value = -1.661232542058096
$ws.Cells.Item(6, 3).Value = 0.001970857450308096
$ws.Cells.Item(7, 2).Value = -0.02163709420448865
$ws.Cells.Item(7, 3).Value = 0.9405622650868235
$ws.Cells.Item(8, 2).Value = 0.001082497821982609
$ws.Cells.Item(8, 3).Value = 0.0000000002857635990099033
$ws.Cells.Item(9, 2).Value = -0.0165653338568605
$ws.Cells.Item(9, 3).Value = 0.01767659614985038
$ws.Cells.Item(10, 2).Value = 0.9565783060681748
$ws.Cells.Item(10, 3).Value = 0.0005658298984455297
$ws.Cells.Item(11, 2).Value = 0.8976101267597042
$ws.Cells.Item(11, 3).Value = 0.002051317654627296
$ws.Cells.Item(12, 2).Value = 0.3413795115798439
$ws.Cells.Item(12, 3).Value = 0.3201842098995435
$ws.Cells.Item(13, 2).Value = -0.00007119155558392746
$ws.Cells.Item(13, 3).Value = 0.5205855714305568
$ws.Cells.Item(14, 2).Value = -0.0000001347001999620896
$ws.Cells.Item(14, 3).Value = 0.1341580450810677
$ws.Cells.Item(15, 2).Value = -0.09502696094831807
$ws.Cells.Item(15, 3).Value = 0.3733135225751852
$ws.Cells.Item(16, 2).Value = 0.12988791604838
$ws.Cells.Item(16, 3).Value = 0.03428659324083805
$ws.Cells.Item(17, 2).Value = -0.04461665532496147
$ws.Cells.Item(17, 3).Value = 0.9952932262971608
$ws.Cells.Item(18, 2).Value = -0.01252667272280466
$ws.Cells.Item(18, 3).Value = 0.1850704120965139
$ws.Cells.Item(19, 2).Value = -0.002282182256763336
$ws.Cells.Item(19, 3).Value = 0.6953222866174456
$ws.Cells.Item(20, 2).Value = 0.4348814034006114
$ws.Cells.Item(20, 3).Value = 0.6790879237786441
$ws.Cells.Item(21, 2).Value = 1.223033574718188
$ws.Cells.Item(21, 3).Value = 0.5648043599060273
$ws.Cells.Item(22, 2).Value = 0.0004818385483290323
$ws.Cells.Item(22, 3).Value = 0.07486298154774534

$ws = $wb.Worksheets.Item("summ0")
$ws.Cells.Item(2, 2).Value = 2.038126794091629
$ws.Cells.Item(2, 3).Value = 0.1193409118109632
$ws.Cells.Item(3, 2).Value = -0.8606579533725857
$ws.Cells.Item(3, 3).Value = 0.2772742808030636
$ws.Cells.Item(4, 2).Value = -1.284220356710692
$ws.Cells.Item(4, 3).Value = 0.006016494247240783
$ws.Cells.Item(5, 2).Value = -1.06041165331014
$ws.Cells.Item(5, 3).Value = 0.02643796602720521
$ws.Cells.Item(6, 2).Value = -1.970688284885907
$ws.Cells.Item(6, 3).Value = 0.0003051578952148187
$ws.Cells.Item(7, 2).Value = 0.2361023748413869
$ws.Cells.Item(7, 3).Value = 0.4981905491929584
$ws.Cells.Item(8, 2).Value = 0.000988272927197344
$ws.Cells.Item(8, 3).Value = 0.000000007216443471525593
$ws.Cells.Item(9, 2).Value = -0.02100637205029261
$ws.Cells.Item(9, 3).Value = 0.002723139226022801
$ws.Cells.Item(10, 2).Value = 0.9106525651853962
$ws.Cells.Item(10, 3).Value = 0.0009681982437368843
$ws.Cells.Item(11, 2).Value = 1.021707395720341
$ws.Cells.Item(11, 3).Value = 0.0004620177505834965
$ws.Cells.Item(12, 2).Value = 0.5262702999445726
$ws.Cells.Item(12, 3).Value = 0.1307320835589011
$ws.Cells.Item(13, 2).Value = -0.0001017766618583697
$ws.Cells.Item(13, 3).Value = 0.3541377952024346
$ws.Cells.Item(14, 2).Value = -0.00000008221893066746802
$ws.Cells.Item(14, 3).Value = 0.3635855773499693
$ws.Cells.Item(15, 2).Value = -0.1009091830523808
$ws.Cells.Item(15, 3).Value = 0.3562772419508039
$ws.Cells.Item(16, 2).Value = 0.1259897472786188
$ws.Cells.Item(16, 3).Value = 0.04381690300395703
$ws.Cells.Item(17, 2).Value = -3.770277722869705
$ws.Cells.Item(17, 3).Value = 0.6159597231375785
$ws.Cells.Item(18, 2).Value = -0.01727060535054269
$ws.Cells.Item(18, 3).Value = 0.07039286759654843
$ws.Cells.Item(19, 2).Value = -0.002744839508495315
$ws.Cells.Item(19, 3).Value = 0.6615442693873745
$ws.Cells.Item(20, 2).Value = 0.8982099795420203
$ws.Cells.Item(20, 3).Value = 0.4043292551333482
$ws.Cells.Item(21, 2).Value = 0.5681895850455498
$ws.Cells.Item(21, 3).Value = 0.7842528926411858
$ws.Cells.Item(22, 2).Value = 0.0005204103870390545
$ws.Cells.Item(22, 3).Value = 0.04958515845897209

# --- Step 2: rename sheet tabs (two-phase to avoid name collisions) ---
$newNames = @("summ14","summ2","summ0","summ1","summ4","summ3","summ5","summ8","summ19")
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = "__tmp_rename_$i"
}
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $newNames[$i-1]
}

Write-Output "done"